$p = $ppt.ActivePresentation

# The deck's slide 21 ("Class-3: Supervised Learning Algorithms") gets speaker
# notes for the first time: a single line with a link to the dataset used for
# the headbrain simple-linear-regression example.
$s = $p.Slides.Item(21)

$notesUrl = "https://www.kaggle.com/codefordata/headbrain-simple-linear-regression"

# Touching the Notes placeholder on the slide's NotesPage is what actually
# materializes ppt/notesSlides/notesSlide2.xml (+ its relationships back to
# slide21 and the notes master) the first time notes are added to this slide.
$notesShape = $s.NotesPage.Shapes.Placeholders.Item(2)
$notesShape.TextFrame.TextRange.Text = $notesUrl

# Make that pasted URL a live hyperlink, matching the authored notes slide.
try {
    $notesRange = $s.NotesPage.Shapes.Item(1).TextFrame.TextRange
    $notesRange.ActionSettings(1).Hyperlink.Address = $notesUrl
} catch {
    # no-op if this host can't attach action settings to a notes placeholder
}
